{"js": "const body = context.document.body;\n\n// --- Change 1: \"Please provide CCTV footage\" gains a trailing sentence and\n// the blank paragraph that used to follow it is folded into it (the author\n// typed across the paragraph break, so the two paragraphs become one). ---\nconst cctvHits = body.search(\"Please provide CCTV footage\", { matchCase: true });\ncctvHits.load(\"text\");\nawait context.sync();\n\nconst cctvRange = cctvHits.items[0];\nconst cctvPara = cctvRange.paragraphs.getFirst();\nconst cctvNextPara = cctvPara.getNext();\ncctvNextPara.load(\"text\");\nawait context.sync();\n\ncctvRange.insertText(\" and a current up to date photograph of prisoner\", \"End\");\ncctvNextPara.delete();\nawait context.sync();\n\n// --- Change 2: \"Please Note, only provide...\" -> \"Please note, only provide...\" ---\nconst noteHits = body.search(\"Please Note, only provide footage for the specific \", {\n  matchCase: true,\n});\nnoteHits.load(\"text\");\nawait context.sync();\nnoteHits.items[0].insertText(\"Please note, only provide footage for the specific \", \"Replace\");\nawait context.sync();\n\n// --- Change 3: drop the manual line break after \"...requested.\" and replace\n// the trailing sentence with the new wording about the photograph. ---\nconst endOfSentence = body.search(\"requested.\", { matchCase: true });\nconst oldTail = body.search(\n  \"Also, please confirm the identity of the individual and the timeframe of when they appear on the footage.\",\n  { matchCase: true }\n);\nendOfSentence.load(\"text\");\noldTail.load(\"text\");\nawait context.sync();\n\nconst lineBreakGap = endOfSentence.items[0]\n  .getRange(\"End\")\n  .expandTo(oldTail.items[0].getRange(\"Start\"));\nlineBreakGap.insertText(\"\", \"Replace\");\nawait context.sync();\n\nconst oldTailAgain = body.search(\n  \"Also, please confirm the identity of the individual and the timeframe of when they appear on the footage.\",\n  { matchCase: true }\n);\noldTailAgain.load(\"text\");\nawait context.sync();\noldTailAgain.items[0].insertText(\n  \" Please also provide a current up-to-date photograph of the prisoner in order to confirm their identity.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1: \"Please provide CCTV footage\" gains a trailing sentence and\n# the blank paragraph that used to follow it is folded into it (the author\n# typed across the paragraph break, so the two paragraphs become one). ---\n$cctv = $d.Content\n$cctv.Find.Text = \"Please provide CCTV footage\"\n$cctv.Find.Execute() | Out-Null\n$cctv.Collapse(0)\n$cctv.InsertAfter(\" and a current up to date photograph of prisoner\")\n\n$allParagraphs = $d.Paragraphs\nfor ($i = 1; $i -le $allParagraphs.Count; $i++) {\n    $p = $allParagraphs.Item($i)\n    if ($p.Range.Text -like \"*Please provide CCTV footage and a current up to date photograph of prisoner*\") {\n        $nextPara = $allParagraphs.Item($i + 1)\n        $nextPara.Range.Delete()\n        break\n    }\n}\n\n# --- Change 2: \"Please Note, only provide...\" -> \"Please note, only provide...\" ---\n$find = $d.Content.Find\n$find.Execute(\n    \"Please Note, only provide footage for the specific \",\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    0,\n    $false,\n    \"Please note, only provide footage for the specific \",\n    2\n) | Out-Null\n\n# --- Change 3: drop the manual line break after \"...requested.\" and replace\n# the trailing sentence with the new wording about the photograph. ---\n$r1 = $d.Content\n$r1.Find.Text = \"requested.\"\n$r1.Find.Execute() | Out-Null\n$endOfR1 = $r1.End\n\n$r2 = $d.Content\n$r2.Find.Text = \"Also, please confirm the identity of the individual and the timeframe of when they appear on the footage.\"\n$r2.Find.Execute() | Out-Null\n$startOfR2 = $r2.Start\n\n$gap = $d.Range($endOfR1, $startOfR2)\n$gap.Text = \"\"\n\n$r3 = $d.Content\n$r3.Find.Text = \"Also, please confirm the identity of the individual and the timeframe of when they appear on the footage.\"\n$r3.Find.Execute() | Out-Null\n$r3.Text = \" Please also provide a current up-to-date photograph of the prisoner in order to confirm their identity.\"\n"}
